$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 37.05583833333333
$ws.Cells.Item(2, 8).Value = 111.167515
$ws.Cells.Item(2, 9).Value = 0.008431126118266585
$ws.Cells.Item(2, 10).Value = 0.008431126118266585
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 17.82379233333333
$ws.Cells.Item(2, 14).Value = 53.471377
$ws.Cells.Item(2, 15).Value = 0.4233776263711466
$ws.Cells.Item(2, 16).Value = 0.4233776263711467
$ws.Cells.Item(2, 17).Value = 660.4755671909061
$ws.Cells.Item(2, 18).Value = 5944.280104718156
$ws.Cells.Item(2, 19).Value = 0.003569550163587486
$ws.Cells.Item(2, 20).Value = 0.003569550163587486

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 37.05583833333333
$ws.Cells.Item(3, 8).Value = 111.167515
$ws.Cells.Item(3, 9).Value = 0.008431126118266585
$ws.Cells.Item(3, 10).Value = 0.008431126118266585
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 10.453073
$ws.Cells.Item(3, 14).Value = 31.359219
$ws.Cells.Item(3, 15).Value = 0.2482971722436279
$ws.Cells.Item(3, 16).Value = 0.2482971722436279
$ws.Cells.Item(3, 17).Value = 387.3473831745317
$ws.Cells.Item(3, 18).Value = 3486.126448570785
$ws.Cells.Item(3, 19).Value = 0.002093424773994988
$ws.Cells.Item(3, 20).Value = 0.002093424773994988

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 37.05583833333333
$ws.Cells.Item(4, 8).Value = 111.167515
$ws.Cells.Item(4, 9).Value = 0.008431126118266585
$ws.Cells.Item(4, 10).Value = 0.008431126118266585
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 4.27602
$ws.Cells.Item(4, 14).Value = 12.82806
$ws.Cells.Item(4, 15).Value = 0.1015704830969034
$ws.Cells.Item(4, 16).Value = 0.1015704830969034
$ws.Cells.Item(4, 17).Value = 158.4515058301
$ws.Cells.Item(4, 18).Value = 1426.0635524709
$ws.Cells.Item(4, 19).Value = 0.0008563535528832572
$ws.Cells.Item(4, 20).Value = 0.0008563535528832573

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 37.05583833333333
$ws.Cells.Item(5, 8).Value = 111.167515
$ws.Cells.Item(5, 9).Value = 0.008431126118266585
$ws.Cells.Item(5, 10).Value = 0.008431126118266585
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 6.530620666666667
$ws.Cells.Item(5, 14).Value = 19.591862
$ws.Cells.Item(5, 15).Value = 0.1551251621919343
$ws.Cells.Item(5, 16).Value = 0.1551251621919343
$ws.Cells.Item(5, 17).Value = 241.9976236403256
$ws.Cells.Item(5, 18).Value = 2177.97861276293
$ws.Cells.Item(5, 19).Value = 0.001307879806556757
$ws.Cells.Item(5, 20).Value = 0.001307879806556757

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 37.05583833333333
$ws.Cells.Item(6, 8).Value = 111.167515
$ws.Cells.Item(6, 9).Value = 0.008431126118266585
$ws.Cells.Item(6, 10).Value = 0.008431126118266585
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 3.015535666666667
$ws.Cells.Item(6, 14).Value = 9.046607
$ws.Cells.Item(6, 15).Value = 0.07162955609638778
$ws.Cells.Item(6, 16).Value = 0.07162955609638778
$ws.Cells.Item(6, 17).Value = 111.7432021524006
$ws.Cells.Item(6, 18).Value = 1005.688819371605
$ws.Cells.Item(6, 19).Value = 0.0006039178212440965
$ws.Cells.Item(6, 20).Value = 0.0006039178212440965

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 50.89916233333333
$ws.Cells.Item(7, 8).Value = 152.697487
$ws.Cells.Item(7, 9).Value = 0.01158082710438721
$ws.Cells.Item(7, 10).Value = 0.01158082710438721
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 17.82379233333333
$ws.Cells.Item(7, 14).Value = 53.471377
$ws.Cells.Item(7, 15).Value = 0.4233776263711466
$ws.Cells.Item(7, 16).Value = 0.4233776263711467
$ws.Cells.Item(7, 17).Value = 907.2160993699554
$ws.Cells.Item(7, 18).Value = 8164.944894329599
$ws.Cells.Item(7, 19).Value = 0.004903063090870098
$ws.Cells.Item(7, 20).Value = 0.004903063090870099

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 50.89916233333333
$ws.Cells.Item(8, 8).Value = 152.697487
$ws.Cells.Item(8, 9).Value = 0.01158082710438721
$ws.Cells.Item(8, 10).Value = 0.01158082710438721
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 10.453073
$ws.Cells.Item(8, 14).Value = 31.359219
$ws.Cells.Item(8, 15).Value = 0.2482971722436279
$ws.Cells.Item(8, 16).Value = 0.2482971722436279
$ws.Cells.Item(8, 17).Value = 532.0526595091836
$ws.Cells.Item(8, 18).Value = 4788.473935582653
$ws.Cells.Item(8, 19).Value = 0.002875486622261707
$ws.Cells.Item(8, 20).Value = 0.002875486622261707

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 50.89916233333333
$ws.Cells.Item(9, 8).Value = 152.697487
$ws.Cells.Item(9, 9).Value = 0.01158082710438721
$ws.Cells.Item(9, 10).Value = 0.01158082710438721
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 4.27602
$ws.Cells.Item(9, 14).Value = 12.82806
$ws.Cells.Item(9, 15).Value = 0.1015704830969034
$ws.Cells.Item(9, 16).Value = 0.1015704830969034
$ws.Cells.Item(9, 17).Value = 217.64583612058
$ws.Cells.Item(9, 18).Value = 1958.81252508522
$ws.Cells.Item(9, 19).Value = 0.001176270203654323
$ws.Cells.Item(9, 20).Value = 0.001176270203654323

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 50.89916233333333
$ws.Cells.Item(10, 8).Value = 152.697487
$ws.Cells.Item(10, 9).Value = 0.01158082710438721
$ws.Cells.Item(10, 10).Value = 0.01158082710438721
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 6.530620666666667
$ws.Cells.Item(10, 14).Value = 19.591862
$ws.Cells.Item(10, 15).Value = 0.1551251621919343
$ws.Cells.Item(10, 16).Value = 0.1551251621919343
$ws.Cells.Item(10, 17).Value = 332.4031214500882
$ws.Cells.Item(10, 18).Value = 2991.628093050794
$ws.Cells.Item(10, 19).Value = 0.001796477682884816
$ws.Cells.Item(10, 20).Value = 0.001796477682884815

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 50.89916233333333
$ws.Cells.Item(11, 8).Value = 152.697487
$ws.Cells.Item(11, 9).Value = 0.01158082710438721
$ws.Cells.Item(11, 10).Value = 0.01158082710438721
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 3.015535666666667
$ws.Cells.Item(11, 14).Value = 9.046607
$ws.Cells.Item(11, 15).Value = 0.07162955609638778
$ws.Cells.Item(11, 16).Value = 0.07162955609638778
$ws.Cells.Item(11, 17).Value = 153.4882394196232
$ws.Cells.Item(11, 18).Value = 1381.394154776609
$ws.Cells.Item(11, 19).Value = 0.0008295295047162721
$ws.Cells.Item(11, 20).Value = 0.0008295295047162721

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 2007.446289
$ws.Cells.Item(12, 8).Value = 6022.338867
$ws.Cells.Item(12, 9).Value = 0.4567440273772037
$ws.Cells.Item(12, 10).Value = 0.4567440273772037
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 17.82379233333333
$ws.Cells.Item(12, 14).Value = 53.471377
$ws.Cells.Item(12, 15).Value = 0.4233776263711466
$ws.Cells.Item(12, 16).Value = 0.4233776263711467
$ws.Cells.Item(12, 17).Value = 35780.30577545665
$ws.Cells.Item(12, 18).Value = 322022.7519791099
$ws.Cells.Item(12, 19).Value = 0.1933752021701585
$ws.Cells.Item(12, 20).Value = 0.1933752021701585

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 2007.446289
$ws.Cells.Item(13, 8).Value = 6022.338867
$ws.Cells.Item(13, 9).Value = 0.4567440273772037
$ws.Cells.Item(13, 10).Value = 0.4567440273772037
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 10.453073
$ws.Cells.Item(13, 14).Value = 31.359219
$ws.Cells.Item(13, 15).Value = 0.2482971722436279
$ws.Cells.Item(13, 16).Value = 0.2482971722436279
$ws.Cells.Item(13, 17).Value = 20983.9826024961
$ws.Cells.Item(13, 18).Value = 188855.8434224649
$ws.Cells.Item(13, 19).Value = 0.1134082504369258
$ws.Cells.Item(13, 20).Value = 0.1134082504369258

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 2007.446289
$ws.Cells.Item(14, 8).Value = 6022.338867
$ws.Cells.Item(14, 9).Value = 0.4567440273772037
$ws.Cells.Item(14, 10).Value = 0.4567440273772037
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 4.27602
$ws.Cells.Item(14, 14).Value = 12.82806
$ws.Cells.Item(14, 15).Value = 0.1015704830969034
$ws.Cells.Item(14, 16).Value = 0.1015704830969034
$ws.Cells.Item(14, 17).Value = 8583.88048068978
$ws.Cells.Item(14, 18).Value = 77254.92432620803
$ws.Cells.Item(14, 19).Value = 0.04639171151232786
$ws.Cells.Item(14, 20).Value = 0.04639171151232787

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 2007.446289
$ws.Cells.Item(15, 8).Value = 6022.338867
$ws.Cells.Item(15, 9).Value = 0.4567440273772037
$ws.Cells.Item(15, 10).Value = 0.4567440273772037
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 6.530620666666667
$ws.Cells.Item(15, 14).Value = 19.591862
$ws.Cells.Item(15, 15).Value = 0.1551251621919343
$ws.Cells.Item(15, 16).Value = 0.1551251621919343
$ws.Cells.Item(15, 17).Value = 13109.87022216671
$ws.Cells.Item(15, 18).Value = 117988.8319995004
$ws.Cells.Item(15, 19).Value = 0.070852491327086
$ws.Cells.Item(15, 20).Value = 0.07085249132708599

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 2007.446289
$ws.Cells.Item(16, 8).Value = 6022.338867
$ws.Cells.Item(16, 9).Value = 0.4567440273772037
$ws.Cells.Item(16, 10).Value = 0.4567440273772037
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 3.015535666666667
$ws.Cells.Item(16, 14).Value = 9.046607
$ws.Cells.Item(16, 15).Value = 0.07162955609638778
$ws.Cells.Item(16, 16).Value = 0.07162955609638778
$ws.Cells.Item(16, 17).Value = 6053.525883397141
$ws.Cells.Item(16, 18).Value = 54481.73295057427
$ws.Cells.Item(16, 19).Value = 0.03271637193070549
$ws.Cells.Item(16, 20).Value = 0.03271637193070549

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 2293.273345666667
$ws.Cells.Item(17, 8).Value = 6879.820037
$ws.Cells.Item(17, 9).Value = 0.5217768014597114
$ws.Cells.Item(17, 10).Value = 0.5217768014597114
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 17.82379233333333
$ws.Cells.Item(17, 14).Value = 53.471377
$ws.Cells.Item(17, 15).Value = 0.4233776263711466
$ws.Cells.Item(17, 16).Value = 0.4233776263711467
$ws.Cells.Item(17, 17).Value = 40874.82787673122
$ws.Cells.Item(17, 18).Value = 367873.450890581
$ws.Cells.Item(17, 19).Value = 0.2209086236975416
$ws.Cells.Item(17, 20).Value = 0.2209086236975417

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 2293.273345666667
$ws.Cells.Item(18, 8).Value = 6879.820037
$ws.Cells.Item(18, 9).Value = 0.5217768014597114
$ws.Cells.Item(18, 10).Value = 0.5217768014597114
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 10.453073
$ws.Cells.Item(18, 14).Value = 31.359219
$ws.Cells.Item(18, 15).Value = 0.2482971722436279
$ws.Cells.Item(18, 16).Value = 0.2482971722436279
$ws.Cells.Item(18, 17).Value = 23971.7536912079
$ws.Cells.Item(18, 18).Value = 215745.7832208711
$ws.Cells.Item(18, 19).Value = 0.1295557043447712
$ws.Cells.Item(18, 20).Value = 0.1295557043447712

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 2293.273345666667
$ws.Cells.Item(19, 8).Value = 6879.820037
$ws.Cells.Item(19, 9).Value = 0.5217768014597114
$ws.Cells.Item(19, 10).Value = 0.5217768014597114
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 4.27602
$ws.Cells.Item(19, 14).Value = 12.82806
$ws.Cells.Item(19, 15).Value = 0.1015704830969034
$ws.Cells.Item(19, 16).Value = 0.1015704830969034
$ws.Cells.Item(19, 17).Value = 9806.082691537582
$ws.Cells.Item(19, 18).Value = 88254.74422383824
$ws.Cells.Item(19, 19).Value = 0.05299712179301996
$ws.Cells.Item(19, 20).Value = 0.05299712179301996

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 2293.273345666667
$ws.Cells.Item(20, 8).Value = 6879.820037
$ws.Cells.Item(20, 9).Value = 0.5217768014597114
$ws.Cells.Item(20, 10).Value = 0.5217768014597114
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 6.530620666666667
$ws.Cells.Item(20, 14).Value = 19.591862
$ws.Cells.Item(20, 15).Value = 0.1551251621919343
$ws.Cells.Item(20, 16).Value = 0.1551251621919343
$ws.Cells.Item(20, 17).Value = 14976.49830552655
$ws.Cells.Item(20, 18).Value = 134788.4847497389
$ws.Cells.Item(20, 19).Value = 0.08094071095442644
$ws.Cells.Item(20, 20).Value = 0.08094071095442643

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 2293.273345666667
$ws.Cells.Item(21, 8).Value = 6879.820037
$ws.Cells.Item(21, 9).Value = 0.5217768014597114
$ws.Cells.Item(21, 10).Value = 0.5217768014597114
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 3.015535666666667
$ws.Cells.Item(21, 14).Value = 9.046607
$ws.Cells.Item(21, 15).Value = 0.07162955609638778
$ws.Cells.Item(21, 16).Value = 0.07162955609638778
$ws.Cells.Item(21, 17).Value = 6915.44756727383
$ws.Cells.Item(21, 18).Value = 62239.02810546446
$ws.Cells.Item(21, 19).Value = 0.03737464066995219
$ws.Cells.Item(21, 20).Value = 0.03737464066995219

$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 7).Value = 6.448603666666666
$ws.Cells.Item(22, 8).Value = 19.345811
$ws.Cells.Item(22, 9).Value = 0.00146721794043115
$ws.Cells.Item(22, 10).Value = 0.00146721794043115
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 17.82379233333333
$ws.Cells.Item(22, 14).Value = 53.471377
$ws.Cells.Item(22, 15).Value = 0.4233776263711466
$ws.Cells.Item(22, 16).Value = 0.4233776263711467
$ws.Cells.Item(22, 17).Value = 114.9385725946385
$ws.Cells.Item(22, 18).Value = 1034.447153351747
$ws.Cells.Item(22, 19).Value = 0.0006211872489889028
$ws.Cells.Item(22, 20).Value = 0.0006211872489889028

$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 7).Value = 6.448603666666666
$ws.Cells.Item(23, 8).Value = 19.345811
$ws.Cells.Item(23, 9).Value = 0.00146721794043115
$ws.Cells.Item(23, 10).Value = 0.00146721794043115
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 10.453073
$ws.Cells.Item(23, 14).Value = 31.359219
$ws.Cells.Item(23, 15).Value = 0.2482971722436279
$ws.Cells.Item(23, 16).Value = 0.2482971722436279
$ws.Cells.Item(23, 17).Value = 67.40772487573433
$ws.Cells.Item(23, 18).Value = 606.6695238816089
$ws.Cells.Item(23, 19).Value = 0.0003643060656741742
$ws.Cells.Item(23, 20).Value = 0.0003643060656741742

$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 7).Value = 6.448603666666666
$ws.Cells.Item(24, 8).Value = 19.345811
$ws.Cells.Item(24, 9).Value = 0.00146721794043115
$ws.Cells.Item(24, 10).Value = 0.00146721794043115
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 4.27602
$ws.Cells.Item(24, 14).Value = 12.82806
$ws.Cells.Item(24, 15).Value = 0.1015704830969034
$ws.Cells.Item(24, 16).Value = 0.1015704830969034
$ws.Cells.Item(24, 17).Value = 27.57435825074
$ws.Cells.Item(24, 18).Value = 248.16922425666
$ws.Cells.Item(24, 19).Value = 0.0001490260350180356
$ws.Cells.Item(24, 20).Value = 0.0001490260350180356

$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 7).Value = 6.448603666666666
$ws.Cells.Item(25, 8).Value = 19.345811
$ws.Cells.Item(25, 9).Value = 0.00146721794043115
$ws.Cells.Item(25, 10).Value = 0.00146721794043115
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 6.530620666666667
$ws.Cells.Item(25, 14).Value = 19.591862
$ws.Cells.Item(25, 15).Value = 0.1551251621919343
$ws.Cells.Item(25, 16).Value = 0.1551251621919343
$ws.Cells.Item(25, 17).Value = 42.11338437667577
$ws.Cells.Item(25, 18).Value = 379.020459390082
$ws.Cells.Item(25, 19).Value = 0.000227602420980298
$ws.Cells.Item(25, 20).Value = 0.0002276024209802979

$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 7).Value = 6.448603666666666
$ws.Cells.Item(26, 8).Value = 19.345811
$ws.Cells.Item(26, 9).Value = 0.00146721794043115
$ws.Cells.Item(26, 10).Value = 0.00146721794043115
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(26, 13).Value = 3.015535666666667
$ws.Cells.Item(26, 14).Value = 9.046607
$ws.Cells.Item(26, 15).Value = 0.07162955609638778
$ws.Cells.Item(26, 16).Value = 0.07162955609638778
$ws.Cells.Item(26, 17).Value = 19.44599435703078
$ws.Cells.Item(26, 18).Value = 175.013949213277
$ws.Cells.Item(26, 19).Value = 0.0001050961697697396
$ws.Cells.Item(26, 20).Value = 0.0001050961697697396
